$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz6")

# Reverse the order of the data rows (rows 2-5), keeping the header row (row 1) fixed.
# Row2 <-> Row5, Row3 <-> Row4

$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$a5 = $ws.Range("A5").Value2
$b5 = $ws.Range("B5").Value2

$ws.Range("A2").Value2 = $a5
$ws.Range("B2").Value2 = $b5
$ws.Range("A3").Value2 = $a4
$ws.Range("B3").Value2 = $b4
$ws.Range("A4").Value2 = $a3
$ws.Range("B4").Value2 = $b3
$ws.Range("A5").Value2 = $a2
$ws.Range("B5").Value2 = $b2

$ws.Range("A6").Select() | Out-Null
